$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157, shifting rows 157:249 down to 158:250.
$ws.Rows.Item(157).Insert()

# Populate the new row 157 with the new data record.
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44438
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = "Fruta"
$ws.Cells.Item(157, 7).Value = 100102
$ws.Cells.Item(157, 8).Value = "Cítricos"
$ws.Cells.Item(157, 9).Value = 100102003
$ws.Cells.Item(157, 10).Value = "Limón"
$ws.Cells.Item(157, 11).Value = "Sin especificar"
$ws.Cells.Item(157, 12).Value = "1a amarillo"
$ws.Cells.Item(157, 13).Value = 300
$ws.Cells.Item(157, 14).Value = 8000
$ws.Cells.Item(157, 15).Value = 8000
$ws.Cells.Item(157, 16).Value = 8000
$ws.Cells.Item(157, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(157, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(157, 19).Value = 500
$ws.Cells.Item(157, 20).Value = 16
